# Realestate Update resale numbers 2023-06-25 08:56
# Adds a new data row (row 75) to the CityResaleNum sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 75

# Force text formatting on Date/Week columns so Excel doesn't
# auto-coerce the string "2023-06-25" into a date serial number,
# or "26" into a plain number (matches existing rows' inline-string types).
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "2023-06-25"
$ws.Cells.Item($row, 2).Value = "08:43:45"
$ws.Cells.Item($row, 3).Value = "Sunday"
$ws.Cells.Item($row, 4).NumberFormat = "@"
$ws.Cells.Item($row, 4).Value = "26"
$ws.Cells.Item($row, 5).Value = 122616
$ws.Cells.Item($row, 6).Value = 134260
$ws.Cells.Item($row, 7).Value = 162813
$ws.Cells.Item($row, 8).Value = 133456
$ws.Cells.Item($row, 9).Value = 177582
$ws.Cells.Item($row, 10).Value = 115835
$ws.Cells.Item($row, 11).Value = 203014
$ws.Cells.Item($row, 12).Value = 225924
$ws.Cells.Item($row, 13).Value = 175853
$ws.Cells.Item($row, 14).Value = 104293
$ws.Cells.Item($row, 15).Value = 39537
$ws.Cells.Item($row, 16).Value = 33799
$ws.Cells.Item($row, 17).Value = 51976
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 35630
$ws.Cells.Item($row, 20).Value = -1
